# Mise en commentaires des paramètres de recherche qu'on n'utilise plus.
# Concretely (per the OOXML diff):
#  - Metadata sheet: insert a new "Jurisdiction" row (with an empty value)
#    right after the "Contact" row, pushing Description/Purpose/Copyright/
#    Immutable down by one row (dimension grows from A1:B14 to A1:B15).
#  - Metadata sheet: bump the "Date" value to the new export timestamp.
#  - "Include from Concepts Medicat" sheet keeps its existing content; its
#    shared-string indices just shift because of the new strings above it.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Insert a new row above the current row 11 ("Description"), so that the
# new row becomes row 11 and everything below shifts down by one.
$ws.Rows.Item(11).Insert()

$ws.Range("A11").Value = "Jurisdiction"
$ws.Range("B11").Value = ""

# Copy the formatting from the row right below (still the plain data style)
# onto the freshly inserted row so it matches the rest of the table instead
# of picking up a default style.
$ws.Range("A12:B12").Copy()
$ws.Range("A11:B11").PasteSpecial(-4122)

# Update the generation Date value (row 8, column B).
$ws.Range("B8").Value = "2024-07-01T07:50:29+00:00"
